$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294, shifting existing rows 294:363 down to 295:364
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new record
$ws.Range("A294").Value = 5
$ws.Range("B294").Value = "Macroferia Regional de Talca"
$ws.Range("C294").Value = "Maule"
$ws.Range("D294").Value = 44754
$ws.Range("E294").Value = 7
$ws.Range("F294").Value = 100114013
$ws.Range("G294").Value = "Zanahoria"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 500
$ws.Range("K294").Value = 9500
$ws.Range("L294").Value = 9500
$ws.Range("M294").Value = 9500
$ws.Range("N294").Value = "$/saco 20 kilos"
$ws.Range("O294").Value = "Región de Ñuble"
$ws.Range("P294").Value = 475
$ws.Range("Q294").Value = 20
$ws.Range("R294").Value = "Hortaliza"
